$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1402
$ws1.Range("F5").Value = 115
$ws1.Range("F6").Value = 69
$ws1.Range("F7").Value = 11878
$ws1.Range("F8").Value = 4436
$ws1.Range("F9").Value = 33
$ws1.Range("F10").Value = 50
$ws1.Range("F12").Value = 20
$ws1.Range("F13").Value = 2565
$ws1.Range("F15").Value = 163
$ws1.Range("F16").Value = 52
$ws1.Range("F17").Value = 5152
$ws1.Range("F19").Value = 193
$ws1.Range("F20").Value = 532
$ws1.Range("F21").Value = 11383
$ws1.Range("F22").Value = 11364
$ws1.Range("F23").Value = 21
$ws1.Range("F28").Value = 23

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1402
$ws4.Range("F5").Value = 115
$ws4.Range("F6").Value = 69
$ws4.Range("F7").Value = 11878
$ws4.Range("F8").Value = 4436
$ws4.Range("F9").Value = 33
$ws4.Range("F10").Value = 50
$ws4.Range("F12").Value = 20
$ws4.Range("F13").Value = 2565
$ws4.Range("F16").Value = 163
$ws4.Range("F17").Value = 52
$ws4.Range("F18").Value = 5152
$ws4.Range("F20").Value = 193
$ws4.Range("F21").Value = 532
$ws4.Range("F22").Value = 11383
$ws4.Range("F23").Value = 11364
$ws4.Range("F24").Value = 21
$ws4.Range("F29").Value = 23
